$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-18: columns B (TOTAL_SUBSTATION_LOAD), C (CONTESTABLE_ENERGY), D (ACTUAL_ENERGY)
$fullRows = @{
  2  = @(78971, 5623.7055, 73347.2945)
  3  = @(75703, 5458.229, 70244.77099999999)
  4  = @(72637, 5507.9745, 67129.0255)
  5  = @(70071, 5475.498, 64595.502)
  6  = @(70844, 5469.52, 65374.48)
  7  = @(73558, 5575.703, 67982.29700000001)
  8  = @(73509, 6141.8735, 67367.1265)
  9  = @(86371, 7163.6355, 79207.3645)
  10 = @(103822, 8647.460500000001, 95174.5395)
  11 = @(94388, 12435.7905, 81952.2095)
  12 = @(97098, 15471.022, 81626.978)
  13 = @(97310, 16384.8475, 80925.1525)
  14 = @(97015, 16062.228, 80952.772)
  15 = @(123772, 16185.204, 107586.796)
  16 = @(122973, 16525.544, 106447.456)
  17 = @(82093, 16462.6175, 65630.38250000001)
  18 = @(56814, 16761.7275, 40052.2725)
}

foreach ($r in $fullRows.Keys) {
  $vals = $fullRows[$r]
  $ws.Cells.Item($r, 2).Value = $vals[0]
  $ws.Cells.Item($r, 3).Value = $vals[1]
  $ws.Cells.Item($r, 4).Value = $vals[2]
}

# Rows 19-25: only column C (CONTESTABLE_ENERGY) changes
$cOnlyRows = @{
  19 = 16532.1625
  20 = 15890.3885
  21 = 13683.775
  22 = 11966.283
  23 = 9398.494000000001
  24 = 6574.9635
  25 = 5831.934499999999
}

foreach ($r in $cOnlyRows.Keys) {
  $ws.Cells.Item($r, 3).Value = $cOnlyRows[$r]
}
